$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.497.38"
$ws.Range("E2").Value = "  +0.91%  "
$ws.Range("D3").Value = "3.446.54"
$ws.Range("E3").Value = "  +2.14%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'578.80"
$ws.Range("E5").Value = "  +1.51%  "
$ws.Range("D6").Value = "'144.19"
$ws.Range("E6").Value = "  +5.64%  "
$ws.Range("D7").Value = "3.447.23"
$ws.Range("E7").Value = "  +2.20%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("E9").Value = "  +1.88%  "
$ws.Range("E10").Value = "  +0.02%  "
$ws.Range("E11").Value = "  +3.22%  "
$ws.Range("E12").Value = "  +2.39%  "
$ws.Range("D13").Value = "4.033.93"
$ws.Range("E13").Value = "  +2.12%  "
$ws.Range("D14").Value = "'27.95"
$ws.Range("E14").Value = "  +9.13%  "
$ws.Range("E15").Value = "  -0.90%  "
$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").Value = "'0.0000173"
$ws.Range("E16").Value = "  +1.65%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "3.449.41"
$ws.Range("E17").Value = "  +2.21%  "
$ws.Range("D18").Value = "61.636.49"
$ws.Range("E18").Value = "  +0.85%  "
$ws.Range("D19").Value = "'6.25"
$ws.Range("E19").Value = "  +9.15%  "
$ws.Range("D20").Value = "'14.27"
$ws.Range("E20").Value = "  +3.64%  "
$ws.Range("D21").Value = "'9.52"
$ws.Range("E21").Value = "  +2.08%  "
$ws.Range("D22").Value = "'389.09"
$ws.Range("E22").Value = "  +4.07%  "
$ws.Range("D23").Value = "'0.564"
$ws.Range("E23").Value = "  +3.25%  "
$ws.Range("D24").Value = "'73.35"
$ws.Range("E24").Value = "  +3.37%  "
$ws.Range("E25").Value = "  +0.15%  "
$ws.Range("E26").Value = "  -0.47%  "
$ws.Range("E27").Value = "  +0.18%  "
$ws.Range("D28").Value = "3.589.81"
$ws.Range("D29").Value = "'0.181"
$ws.Range("E29").Value = "  +1.60%  "
$ws.Range("E30").Value = "  +3.52%  "
$ws.Range("E31").Value = "  +0.06%  "
$ws.Range("E32").Value = "  +1.50%  "
$ws.Range("E33").Value = "  +2.72%  "
$ws.Range("E34").Value = "  -11.14%  "
$ws.Range("E35").Value = "  -0.05%  "
$ws.Range("D36").Value = "'24.04"
$ws.Range("E36").Value = "  +3.31%  "
$ws.Range("D37").Value = "3.475.65"
$ws.Range("E37").Value = "  +2.30%  "
$ws.Range("D38").Value = "'7.00"
$ws.Range("E38").Value = "  +3.45%  "
$ws.Range("D39").Value = "'5.12"
$ws.Range("E39").Value = "  +0.22%  "
$ws.Range("E40").Value = "  +0.69%  "
$ws.Range("D41").Value = "'166.70"
$ws.Range("E41").Value = "  +1.13%  "
$ws.Range("D42").Value = "'28.12"
$ws.Range("E42").Value = "  +13.46%  "
$ws.Range("D43").Value = "'0.0782"
$ws.Range("E43").Value = "  +3.33%  "
$ws.Range("E44").Value = "  +3.79%  "
$ws.Range("E45").Value = "  +0.01%  "
$ws.Range("E46").Value = "  +1.55%  "
$ws.Range("D47").Value = "'4.48"
$ws.Range("E47").Value = "  +3.97%  "
$ws.Range("E48").Value = "  +2.71%  "
$ws.Range("D49").Value = "2.581.14"
$ws.Range("E49").Value = "  +1.51%  "
$ws.Range("E50").Value = "  -1.72%  "
$ws.Range("E51").Value = "  +2.26%  "
